# Auto-generated edit script applying numeric corrections to the
# Jenova_Profits workbook's Leve-profit calculation columns (H-N) across
# all 8 job sheets, per the scheduled runner's price-refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 50939.824
$ws.Range("J17").Value = 52238.793
$ws.Range("L17").Value = 156716.379
$ws.Range("N17").Value = -157052.379

$ws.Range("H64").Value = 6709.778

$ws.Range("H67").Value = 6709.778

$ws.Range("H86").Value = 3098294.5
$ws.Range("I86").Value = 1919.7778
$ws.Range("J86").Value = 6581716
$ws.Range("K86").Value = 1919.7778
$ws.Range("L86").Value = 6581716
$ws.Range("M86").Value = -796.7778000000001
$ws.Range("N86").Value = -6583962

$ws.Range("H89").Value = 3098294.5
$ws.Range("I89").Value = 1919.7778
$ws.Range("J89").Value = 6581716
$ws.Range("K89").Value = 9598.889000000001
$ws.Range("L89").Value = 32908580
$ws.Range("M89").Value = -3982.889000000001
$ws.Range("N89").Value = -32919812

$ws.Range("H106").Value = 2105.3845
$ws.Range("I106").Value = 2173
$ws.Range("K106").Value = 2173
$ws.Range("M106").Value = -1542

$ws.Range("H125").Value = 7891.3076
$ws.Range("I125").Value = 7517.6665
$ws.Range("K125").Value = 67658.9985
$ws.Range("M125").Value = -65198.9985

$ws.Range("H137").Value = 3155.606
$ws.Range("I137").Value = 2513.1538
$ws.Range("J137").Value = 3573.2
$ws.Range("K137").Value = 7539.4614
$ws.Range("L137").Value = 10719.6
$ws.Range("M137").Value = -4989.4614
$ws.Range("N137").Value = -15819.6

$ws.Range("H141").Value = 1992.6875
$ws.Range("I141").Value = 1992.6875
$ws.Range("K141").Value = 5978.0625
$ws.Range("M141").Value = -798.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 6994
$ws.Range("J10").Value = 4000
$ws.Range("L10").Value = 4000
$ws.Range("N10").Value = -4340

$ws.Range("H61").Value = 1105.909
$ws.Range("I61").Value = 1129.5555
$ws.Range("K61").Value = 1129.5555
$ws.Range("M61").Value = -917.5554999999999

$ws.Range("H110").Value = 274219.7
$ws.Range("I110").Value = 376620.44
$ws.Range("J110").Value = 1151
$ws.Range("K110").Value = 376620.44
$ws.Range("L110").Value = 1151
$ws.Range("M110").Value = -374575.44
$ws.Range("N110").Value = -5241

$ws.Range("H132").Value = 6144.1626
$ws.Range("I132").Value = 4062.8286
$ws.Range("J132").Value = 15250
$ws.Range("K132").Value = 12188.4858
$ws.Range("L132").Value = 45750
$ws.Range("M132").Value = -9658.485799999999
$ws.Range("N132").Value = -50810

$ws.Range("H136").Value = 1105.909
$ws.Range("I136").Value = 1129.5555
$ws.Range("K136").Value = 3388.6665
$ws.Range("M136").Value = -838.6664999999998

$ws.Range("H138").Value = 62500
$ws.Range("J138").Value = 62500
$ws.Range("L138").Value = 62500
$ws.Range("N138").Value = -72780

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 63186.4
$ws.Range("J74").Value = 60239.445
$ws.Range("L74").Value = 60239.445
$ws.Range("N74").Value = -62111.445

$ws.Range("H77").Value = 63186.4
$ws.Range("J77").Value = 60239.445
$ws.Range("L77").Value = 180718.335
$ws.Range("N77").Value = -190078.335

$ws.Range("H105").Value = 45720.043
$ws.Range("I105").Value = 252554.25
$ws.Range("K105").Value = 252554.25
$ws.Range("M105").Value = -250807.25

$ws.Range("H107").Value = 438.86957
$ws.Range("I107").Value = 407.57144
$ws.Range("K107").Value = 407.57144
$ws.Range("M107").Value = 1512.42856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2914.8064
$ws.Range("I31").Value = 2175.5334
$ws.Range("K31").Value = 2175.5334
$ws.Range("M31").Value = -1880.5334

$ws.Range("H34").Value = 2914.8064
$ws.Range("I34").Value = 2175.5334
$ws.Range("K34").Value = 2175.5334
$ws.Range("M34").Value = -1973.5334

$ws.Range("H48").Value = 2500
$ws.Range("J48").Value = 2500
$ws.Range("L48").Value = 2500
$ws.Range("N48").Value = -3452

$ws.Range("H55").Value = 17632.2
$ws.Range("I55").Value = 14500
$ws.Range("J55").Value = 19720.334
$ws.Range("K55").Value = 14500
$ws.Range("L55").Value = 19720.334
$ws.Range("M55").Value = -14185
$ws.Range("N55").Value = -20350.334

$ws.Range("H57").Value = 8000
$ws.Range("J57").Value = 8000
$ws.Range("L57").Value = 8000
$ws.Range("N57").Value = -9120

$ws.Range("H105").Value = 1278.2174
$ws.Range("I105").Value = 1209.95
$ws.Range("K105").Value = 1209.95
$ws.Range("M105").Value = 537.05

$ws.Range("H134").Value = 5165
$ws.Range("I134").Value = 4399.7646
$ws.Range("J134").Value = 7333.1665
$ws.Range("K134").Value = 13199.2938
$ws.Range("L134").Value = 21999.4995
$ws.Range("M134").Value = -10664.2938
$ws.Range("N134").Value = -27069.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1400034.2
$ws.Range("I6").Value = 49.5
$ws.Range("J6").Value = 2333357.2
$ws.Range("K6").Value = 148.5
$ws.Range("L6").Value = 7000071.600000001
$ws.Range("M6").Value = -35.5
$ws.Range("N6").Value = -7000297.600000001

$ws.Range("H86").Value = 3666698.2
$ws.Range("I86").Value = 92
$ws.Range("K86").Value = 276
$ws.Range("M86").Value = 910

$ws.Range("H89").Value = 3666698.2
$ws.Range("I89").Value = 92
$ws.Range("K89").Value = 828
$ws.Range("M89").Value = 5100

$ws.Range("H122").Value = 41074.32
$ws.Range("J122").Value = 56642.11
$ws.Range("L122").Value = 509778.99
$ws.Range("N122").Value = -514678.99

$ws.Range("H131").Value = 4688.0557
$ws.Range("I131").Value = 1469.5834
$ws.Range("J131").Value = 11125
$ws.Range("K131").Value = 4408.7502
$ws.Range("L131").Value = 33375
$ws.Range("M131").Value = 631.2497999999996
$ws.Range("N131").Value = -43455

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 140.4
$ws.Range("I2").Value = 140.4
$ws.Range("K2").Value = 140.4
$ws.Range("M2").Value = -27.40000000000001

$ws.Range("H102").Value = 3521.9565
$ws.Range("I102").Value = 3940.7646
$ws.Range("J102").Value = 2335.3333
$ws.Range("K102").Value = 3940.7646
$ws.Range("L102").Value = 2335.3333
$ws.Range("M102").Value = -2318.7646
$ws.Range("N102").Value = -5579.3333

$ws.Range("H122").Value = 9343.929
$ws.Range("I122").Value = 13203
$ws.Range("J122").Value = 7200
$ws.Range("K122").Value = 39609
$ws.Range("L122").Value = 21600
$ws.Range("M122").Value = -37159
$ws.Range("N122").Value = -26500

$ws.Range("H132").Value = 420473.1
$ws.Range("I132").Value = 529439.7
$ws.Range("K132").Value = 1588319.1
$ws.Range("M132").Value = -1585789.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 78763.64
$ws.Range("J68").Value = 251197.75
$ws.Range("L68").Value = 251197.75
$ws.Range("N68").Value = -252695.75

$ws.Range("H71").Value = 78763.64
$ws.Range("J71").Value = 251197.75
$ws.Range("L71").Value = 1255988.75
$ws.Range("N71").Value = -1263476.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21365.334
$ws.Range("J54").Value = 21365.334
$ws.Range("L54").Value = 21365.334
$ws.Range("N54").Value = -22405.334

$ws.Range("H81").Value = 3929.1428
$ws.Range("I81").Value = 1521
$ws.Range("J81").Value = 9949.5
$ws.Range("K81").Value = 3042
$ws.Range("L81").Value = 19899
$ws.Range("M81").Value = -1981
$ws.Range("N81").Value = -22021

$ws.Range("H84").Value = 3929.1428
$ws.Range("I84").Value = 1521
$ws.Range("J84").Value = 9949.5
$ws.Range("K84").Value = 15210
$ws.Range("L84").Value = 99495
$ws.Range("M84").Value = -9906
$ws.Range("N84").Value = -110103

$ws.Range("H96").Value = 101527.6
$ws.Range("J96").Value = 1500
$ws.Range("L96").Value = 1500
$ws.Range("N96").Value = -4246

$ws.Range("H107").Value = 58168.777
$ws.Range("I107").Value = 80053.30499999999
$ws.Range("K107").Value = 240159.915
$ws.Range("M107").Value = -238239.915

